$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "last updated" timestamp string (A1)
$ws.Range("A1").Value = "Datos actualizados a 23 de Julio de 2020 a las 21:04"

# Row 4: Estados Unidos
$ws.Range("A4").Value = "Estados Unidos"
$ws.Range("B4").Value = 4141163
$ws.Range("C4").Value = 40288
$ws.Range("D4").Value = 1954790
$ws.Range("E4").Value = 2039603
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 587
$ws.Range("H4").Value = 146770

# Row 6: India
$ws.Range("A6").Value = "India"
$ws.Range("B6").Value = 1287083
$ws.Range("C6").Value = 47399
$ws.Range("D6").Value = 816205
$ws.Range("E6").Value = 440233
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 755
$ws.Range("H6").Value = 30645

# Row 21: Alemania
$ws.Range("A21").Value = "Alemania"
$ws.Range("B21").Value = 204889
$ws.Range("C21").Value = 419
$ws.Range("D21").Value = 189000
$ws.Range("E21").Value = 6703
$ws.Range("F21").Value = 0
$ws.Range("G21").Value = 4
$ws.Range("H21").Value = 9186

# Row 22: Francia
$ws.Range("A22").Value = "Francia"
$ws.Range("B22").Value = 179398
$ws.Range("C22").Value = 1062
$ws.Range("D22").Value = 80472
$ws.Range("E22").Value = 68744
$ws.Range("F22").Value = 0
$ws.Range("G22").Value = 10
$ws.Range("H22").Value = 30182

# Row 65: Uzbekistan
$ws.Range("A65").Value = "Uzbekistan"
$ws.Range("B65").Value = 18703
$ws.Range("C65").Value = 324
$ws.Range("D65").Value = 10002
$ws.Range("E65").Value = 8599
$ws.Range("F65").Value = 0
$ws.Range("G65").Value = 4
$ws.Range("H65").Value = 102

# Row 91: Guayana Francesa
$ws.Range("A91").Value = "Guayana Francesa"
$ws.Range("B91").Value = 7086
$ws.Range("C91").Value = 203
$ws.Range("D91").Value = 5376
$ws.Range("E91").Value = 1670
$ws.Range("F91").Value = 0
$ws.Range("G91").Value = 1
$ws.Range("H91").Value = 40

# Row 92: Tayikistan
$ws.Range("A92").Value = "Tayikistan"
$ws.Range("B92").Value = 7060
$ws.Range("C92").Value = 45
$ws.Range("D92").Value = 5793
$ws.Range("E92").Value = 1209
$ws.Range("F92").Value = 0
$ws.Range("G92").Value = 0
$ws.Range("H92").Value = 58

# Row 97: Republica de Yibuti
$ws.Range("A97").Value = "Republica de Yibuti"
$ws.Range("B97").Value = 5031
$ws.Range("C97").Value = 1
$ws.Range("D97").Value = 4927
$ws.Range("E97").Value = 46
$ws.Range("F97").Value = 0
$ws.Range("G97").Value = 0
$ws.Range("H97").Value = 58

# Row 118: Libia
$ws.Range("A118").Value = "Libia"
$ws.Range("B118").Value = 2314
$ws.Range("C118").Value = 138
$ws.Range("D118").Value = 501
$ws.Range("E118").Value = 1757
$ws.Range("F118").Value = 0
$ws.Range("G118").Value = 3
$ws.Range("H118").Value = 56

# Row 119: Hong Kong
$ws.Range("A119").Value = "Hong Kong"
$ws.Range("B119").Value = 2250
$ws.Range("C119").Value = 118
$ws.Range("D119").Value = 1379
$ws.Range("E119").Value = 856
$ws.Range("F119").Value = 0
$ws.Range("G119").Value = 1
$ws.Range("H119").Value = 15

# Row 120: Sudan del Sur
$ws.Range("A120").Value = "Sudan del Sur"
$ws.Range("B120").Value = 2239
$ws.Range("C120").Value = 28
$ws.Range("D120").Value = 1175
$ws.Range("E120").Value = 1019
$ws.Range("F120").Value = 0
$ws.Range("G120").Value = 0
$ws.Range("H120").Value = 45

# Row 126: Suazilandia
$ws.Range("A126").Value = "Suazilandia"
$ws.Range("B126").Value = 2021
$ws.Range("C126").Value = 83
$ws.Range("D126").Value = 882
$ws.Range("E126").Value = 1111
$ws.Range("F126").Value = 0
$ws.Range("G126").Value = 3
$ws.Range("H126").Value = 28

# Row 127: Lituania
$ws.Range("A127").Value = "Lituania"
$ws.Range("B127").Value = 1960
$ws.Range("C127").Value = 9
$ws.Range("D127").Value = 1611
$ws.Range("E127").Value = 269
$ws.Range("F127").Value = 0
$ws.Range("G127").Value = 0
$ws.Range("H127").Value = 80

# Row 128: Guinea-Bisau
$ws.Range("A128").Value = "Guinea-Bisau"
$ws.Range("B128").Value = 1954
$ws.Range("C128").Value = 0
$ws.Range("D128").Value = 803
$ws.Range("E128").Value = 1125
$ws.Range("F128").Value = 0
$ws.Range("G128").Value = 0
$ws.Range("H128").Value = 26

# Row 150: Angola
$ws.Range("A150").Value = "Angola"
$ws.Range("B150").Value = 851
$ws.Range("C150").Value = 39
$ws.Range("D150").Value = 236
$ws.Range("E150").Value = 582
$ws.Range("F150").Value = 0
$ws.Range("G150").Value = 0
$ws.Range("H150").Value = 33

# Row 151: Jamaica
$ws.Range("A151").Value = "Jamaica"
$ws.Range("B151").Value = 816
$ws.Range("C151").Value = 6
$ws.Range("D151").Value = 710
$ws.Range("E151").Value = 96
$ws.Range("F151").Value = 0
$ws.Range("G151").Value = 0
$ws.Range("H151").Value = 10

# Row 168: Comoras
$ws.Range("A168").Value = "Comoras"
$ws.Range("B168").Value = 340
$ws.Range("C168").Value = 3
$ws.Range("D168").Value = 324
$ws.Range("E168").Value = 9
$ws.Range("F168").Value = 0
$ws.Range("G168").Value = 0
$ws.Range("H168").Value = 7

# Row 172: Eritrea
$ws.Range("A172").Value = "Eritrea"
$ws.Range("B172").Value = 261
$ws.Range("C172").Value = 10
$ws.Range("D172").Value = 189
$ws.Range("E172").Value = 72
$ws.Range("F172").Value = 0
$ws.Range("G172").Value = 0
$ws.Range("H172").Value = 0

# Row 210: Groenlandia
$ws.Range("A210").Value = "Groenlandia"
$ws.Range("B210").Value = 13
$ws.Range("C210").Value = 0
$ws.Range("D210").Value = 13
$ws.Range("E210").Value = 0
$ws.Range("F210").Value = 0
$ws.Range("G210").Value = 0
$ws.Range("H210").Value = 0

# Row 211: Islas Malvinas
$ws.Range("A211").Value = "Islas Malvinas"
$ws.Range("B211").Value = 13
$ws.Range("C211").Value = 0
$ws.Range("D211").Value = 13
$ws.Range("E211").Value = 0
$ws.Range("F211").Value = 0
$ws.Range("G211").Value = 0
$ws.Range("H211").Value = 0
